# Correction in TSR formula
#
# The TSR column (F) on Sheet1 used the formula  SUM(1, -E, -D)  which is
# being replaced with the corrected weighted formula
#   (5 * SUM(1, -E) + 4 * SUM(1, -D)) / 9
# Also fixes a wrong input value in E10 and tidies a couple of cell
# styles that were left inconsistent by the previous formula range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Data correction -------------------------------------------------
# E10 should be 1 (was 0.75)
$ws.Range("E10").Value = 1

# --- Formula correction ------------------------------------------------
# F4 holds its own (non-shared) copy of the corrected formula.
$ws.Range("F4").Formula = "=(5 * SUM(1, -E4) + 4 * SUM(1, -D4)) / 9"

# F5:F36 all share one formula. Fill the whole block first so the engine
# builds a single shared-formula group covering F5:F36...
$ws.Range("F5:F36").Formula = "=(5 * SUM(1, -E5) + 4 * SUM(1, -D5)) / 9"

# ...then clear back out the separator rows that must stay empty (F9,
# F20, F26 and F30 have no formula in the corrected sheet; F28 keeps
# its formula, unlike before).
$ws.Range("F9").ClearContents() | Out-Null
$ws.Range("F20").ClearContents() | Out-Null
$ws.Range("F26").ClearContents() | Out-Null
$ws.Range("F30").ClearContents() | Out-Null

# --- Style touch-ups -----------------------------------------------------
# F7 previously had no explicit style; it now matches the rest of the
# column (same formatting as F6/F8).
$ws.Range("F6").Copy() | Out-Null
$ws.Range("F7").PasteSpecial($xlPasteFormats) | Out-Null

# F26 and F30 move to the same (border-less-on-top) style used by the
# rest of those separator rows, e.g. E26/E30.
$ws.Range("E26").Copy() | Out-Null
$ws.Range("F26").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("E30").Copy() | Out-Null
$ws.Range("F30").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Application.CutCopyMode = 0

# --- Selection -----------------------------------------------------------
# The workbook was left with F4 selected.
$ws.Range("F4").Select() | Out-Null
